# Daily attendance processing - 2025-10-10 23:40:35
#
# For every "Recorded By" (column G) cell that lists multiple recorders
# separated by ", ", the first and last recorder in the list are swapped.
# Cells with only a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 1) {
            $first = $parts[0]
            $last = $parts[$parts.Length - 1]
            $parts[0] = $last
            $parts[$parts.Length - 1] = $first

            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
